$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string driven country ordering (swap pairs of countries) ---
# Rows 111/112: Zimbabue <-> Mozambique swap display order
$ws.Range("A111").Value = "Mozambique"
$ws.Range("A112").Value = "Zimbabue"

# Rows 184/185: Isla de Man <-> Curazao swap display order
$ws.Range("A184").Value = "Curazao"
$ws.Range("A185").Value = "Isla de Man"

# Rows 207/208: Santa Lucia <-> Timor Oriental swap display order (data identical, no value changes)
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Update statistic values (new daily COVID case numbers) ---
$ws.Range("B4").Value = 7295718
$ws.Range("C4").Value = 8157
$ws.Range("D4").Value = 4525753
$ws.Range("E4").Value = 2560723
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = 209242

$ws.Range("B6").Value = 4719099
$ws.Range("C6").Value = 984
$ws.Range("E6").Value = 526759
$ws.Range("G6").Value = 62
$ws.Range("H6").Value = 141503

$ws.Range("B17").Value = 434969
$ws.Range("C17").Value = 5693
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 41988

$ws.Range("B21").Value = 314433
$ws.Range("C21").Value = 1467
$ws.Range("D21").Value = 275630
$ws.Range("E21").Value = 30806
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = 7997

$ws.Range("B25").Value = 285729
$ws.Range("C25").Value = 704
$ws.Range("E25").Value = 26695
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 9534

$ws.Range("B29").Value = 153058
$ws.Range("C29").Value = 1387
$ws.Range("D29").Value = 131086
$ws.Range("E29").Value = 12704
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 9268

$ws.Range("B57").Value = 64012
$ws.Range("C57").Value = 718
$ws.Range("D57").Value = 31220
$ws.Range("E57").Value = 32191
$ws.Range("G57").Value = 10
$ws.Range("H57").Value = 601

$ws.Range("D61").Value = 42700
$ws.Range("E61").Value = 7100

$ws.Range("B62").Value = 51067
$ws.Range("C62").Value = 153
$ws.Range("D62").Value = 35860
$ws.Range("E62").Value = 13493
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 1714

$ws.Range("B65").Value = 46387
$ws.Range("C65").Value = 165
$ws.Range("D65").Value = 45618
$ws.Range("E65").Value = 470

$ws.Range("B73").Value = 36254
$ws.Range("C73").Value = 1012
$ws.Range("D73").Value = 16089
$ws.Range("E73").Value = 19818
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 347

$ws.Range("B74").Value = 34990
$ws.Range("C74").Value = 430
$ws.Range("E74").Value = 9824

$ws.Range("B87").Value = 17444
$ws.Range("C87").Value = 216
$ws.Range("E87").Value = 7076
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 379

$ws.Range("B101").Value = 10313
$ws.Range("C101").Value = 116
$ws.Range("D101").Value = 6456
$ws.Range("E101").Value = 3699

$ws.Range("B110").Value = 8357
$ws.Range("C110").Value = 46
$ws.Range("E110").Value = 1257

$ws.Range("B111").Value = 7983
$ws.Range("C111").Value = 226
$ws.Range("D111").Value = 4807
$ws.Range("E111").Value = 3118
$ws.Range("G111").Value = 4
$ws.Range("H111").Value = 58

$ws.Range("B112").Value = 7803
$ws.Range("D112").Value = 6067
$ws.Range("E112").Value = 1509
$ws.Range("H112").Value = 227

$ws.Range("B148").Value = 2686
$ws.Range("C148").Value = 10
$ws.Range("E148").Value = 1347

$ws.Range("B160").Value = 1696
$ws.Range("C160").Value = 12
$ws.Range("E160").Value = 305

$ws.Range("B164").Value = 1339
$ws.Range("C164").Value = 1
$ws.Range("E164").Value = 36

$ws.Range("B184").Value = 360
$ws.Range("C184").Value = 23
$ws.Range("D184").Value = 141
$ws.Range("E184").Value = 218
$ws.Range("H184").Value = 1

$ws.Range("B185").Value = 340
$ws.Range("D185").Value = 314
$ws.Range("E185").Value = 2
$ws.Range("H185").Value = 24

$ws.Range("B196").Value = 101
$ws.Range("C196").Value = 3
$ws.Range("E196").Value = 6
